# Add a new PROFIT column (K) to the SalesData sheet, computed from the
# PRODUCTLINE (E) and SALES (H) columns, and extend the sheet's dimension /
# autofilter / filter-database defined name to cover the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SalesData")

# Header
$ws.Range("K1").Value = "PROFIT"

# Fill the PROFIT formula down for every data row (2..107), each row
# referencing its own E/H cells (mirrors a fill-down of the K2 formula).
for ($r = 2; $r -le 107; $r++) {
    $formula = '=IF(E' + $r + '="INTERNALSTORAGE",H' + $r + '*52%,' +
               'IF(E' + $r + '="SERVICE",H' + $r + '*100%,' +
               'IF(E' + $r + '="DISPLAY",H' + $r + '*50%,' +
               'IF(E' + $r + '="RAM",H' + $r + '*65%,' +
               'IF(E' + $r + '="SOFTWARE",H' + $r + '*100%,' +
               'IF(E' + $r + '="BATTERY",H' + $r + '*70%,' +
               'IF(E' + $r + '="KEYBOARD",H' + $r + '*60%,' +
               'IF(E' + $r + '="ADAPTOR",H' + $r + '*45%,' +
               'IF(E' + $r + '="MAINBOARD",H' + $r + '*15%,' +
               'IF(E' + $r + '="ACCESSORIES",H' + $r + '*25%,' +
               'IF(E' + $r + '="VGA",H' + $r + '*15%,' +
               'IF(E' + $r + '="POWERSUPPLY",H' + $r + '*40%,' +
               'IF(E' + $r + '="PROCESSOR",H' + $r + '*20%,' +
               'IF(E' + $r + '="SECOND",H' + $r + '*150%,' +
               'IF(E' + $r + '="CASING",H' + $r + '*15%,' +
               'IF(E' + $r + '="MONITOR",H' + $r + '*5%,' +
               'IF(E' + $r + '="EXTERNALSTORAGE",100000,H' + $r + '))))))))))))))))'
    $ws.Range("K$r").Formula = $formula
}

# Extend the AutoFilter to include the new column.
$ws.AutoFilterMode = $false
$ws.Range("A1:K107").AutoFilter(1)

# Keep the (hidden) _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "SalesData!_FilterDatabase") {
        $n.RefersTo = "=SalesData!`$A`$1:`$K`$107"
    }
}

# Match the author's final selection / scroll position.
$ws.Range("K2:K107").Select()
